# "testes -> teste de cadastro de empresas"
# Adds a new row (12) to Planilha1 recording a new "empresa" test entry,
# following the same pattern as the existing rows (LOGIN, QTDD FUNC, PRAZO,
# DT INICIO, NOTA).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$xlPasteFormats = -4122

# Set E12 before A12 so new shared strings are inserted in the same order
# as in the target workbook ("Teste10 - Gui" then "EmpT10").
$ws.Range("E12").Value = "Teste10 - Gui"
$ws.Range("A12").Value = "EmpT10"
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 2

# Copy the date formatting from the cell above (D11) instead of assigning a
# NumberFormat string directly, so the existing date style is reused rather
# than a brand new style being created.
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("D12").Value = (Get-Date -Year 2025 -Month 2 -Day 27 -Hour 0 -Minute 0 -Second 0).Date

$ws.Range("A12").Select()
